$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Value)
    # Force the cell to stay text (matches source t="inlineStr" cells),
    # even when the new value looks numeric (e.g. "1.01"), then restore
    # the default "Normal" style so no stray numFmt/style index sticks.
    $Cell.NumberFormat = "@"
    $Cell.Value = $Value
    $Cell.Style = "Normal"
}

Set-TextValue $ws.Cells.Item(2, 4) '26.698.33'
$ws.Cells.Item(2, 5).Value = '  +0.01%  '
Set-TextValue $ws.Cells.Item(3, 4) '1.647.20'
$ws.Cells.Item(3, 5).Value = '  +0.66%  '
Set-TextValue $ws.Cells.Item(4, 4) '1.01'
$ws.Cells.Item(4, 5).Value = '  +0.36%  '
Set-TextValue $ws.Cells.Item(5, 4) '216.19'
$ws.Cells.Item(5, 5).Value = '  +1.23%  '
$ws.Cells.Item(6, 5).Value = '  -0.54%  '
Set-TextValue $ws.Cells.Item(7, 4) '1.01'
$ws.Cells.Item(7, 5).Value = '  +0.34%  '
$ws.Cells.Item(8, 5).Value = '  -0.43%  '
$ws.Cells.Item(9, 5).Value = '  +0.50%  '
Set-TextValue $ws.Cells.Item(10, 4) '19.41'
$ws.Cells.Item(10, 5).Value = '  +0.41%  '
$ws.Cells.Item(11, 5).Value = '  +0.14%  '
Set-TextValue $ws.Cells.Item(12, 4) '1.877.00'
$ws.Cells.Item(12, 5).Value = '  +0.70%  '
$ws.Cells.Item(13, 2).Value = 'Polkadot'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Cells.Item(13, 4) '4.24'
$ws.Cells.Item(13, 5).Value = '  +3.20%  '
$ws.Cells.Item(14, 2).Value = 'WrappedEther'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Cells.Item(14, 4) '1.626.56'
$ws.Cells.Item(14, 5).Value = '  -0.75%  '
Set-TextValue $ws.Cells.Item(15, 4) '0.536'
$ws.Cells.Item(15, 5).Value = '  +1.46%  '
Set-TextValue $ws.Cells.Item(16, 4) '66.41'
$ws.Cells.Item(16, 5).Value = '  +4.55%  '
Set-TextValue $ws.Cells.Item(17, 4) '26.769.53'
$ws.Cells.Item(17, 5).Value = '  +0.33%  '
$ws.Cells.Item(18, 5).Value = '  +1.27%  '
Set-TextValue $ws.Cells.Item(19, 4) '220.13'
$ws.Cells.Item(19, 5).Value = '  +0.09%  '
$ws.Cells.Item(20, 5).Value = '  +0.30%  '
Set-TextValue $ws.Cells.Item(21, 4) '4.41'
$ws.Cells.Item(21, 5).Value = '  +1.91%  '
Set-TextValue $ws.Cells.Item(22, 4) '6.34'
$ws.Cells.Item(22, 5).Value = '  +2.25%  '
Set-TextValue $ws.Cells.Item(23, 4) '9.59'
$ws.Cells.Item(23, 5).Value = '  +1.40%  '
Set-TextValue $ws.Cells.Item(24, 4) '2.11'
$ws.Cells.Item(24, 5).Value = '  +9.88%  '
Set-TextValue $ws.Cells.Item(25, 4) '147.12'
$ws.Cells.Item(25, 5).Value = '  -0.71%  '
Set-TextValue $ws.Cells.Item(26, 4) '1.01'
$ws.Cells.Item(26, 5).Value = '  +0.39%  '
$ws.Cells.Item(27, 5).Value = '  -0.75%  '
Set-TextValue $ws.Cells.Item(28, 4) '7.12'
$ws.Cells.Item(28, 5).Value = '  +2.52%  '
Set-TextValue $ws.Cells.Item(29, 4) '15.91'
$ws.Cells.Item(29, 5).Value = '  +2.36%  '
Set-TextValue $ws.Cells.Item(30, 4) '0.0519'
$ws.Cells.Item(30, 5).Value = '  +1.67%  '
$ws.Cells.Item(31, 5).Value = '  +0.74%  '
Set-TextValue $ws.Cells.Item(32, 4) '3.41'
$ws.Cells.Item(32, 5).Value = '  +2.19%  '
Set-TextValue $ws.Cells.Item(33, 4) '3.07'
$ws.Cells.Item(33, 5).Value = '  +2.64%  '
Set-TextValue $ws.Cells.Item(34, 4) '1.288.97'
$ws.Cells.Item(34, 5).Value = '  +6.17%  '
Set-TextValue $ws.Cells.Item(35, 4) '1.55'
$ws.Cells.Item(35, 5).Value = '  +1.85%  '
$ws.Cells.Item(36, 5).Value = '  +6.29%  '
$ws.Cells.Item(37, 5).Value = '  +0.53%  '
$ws.Cells.Item(38, 2).Value = 'ImmutableX'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Cells.Item(38, 4) '0.529'
$ws.Cells.Item(38, 5).Value = '  +4.26%  '
$ws.Cells.Item(39, 2).Value = 'ARBITRUM'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Cells.Item(39, 4) '0.829'
$ws.Cells.Item(39, 5).Value = '  +2.08%  '
Set-TextValue $ws.Cells.Item(40, 4) '1.01'
$ws.Cells.Item(40, 5).Value = '  +0.35%  '
Set-TextValue $ws.Cells.Item(41, 4) '0.813'
$ws.Cells.Item(41, 5).Value = '  +2.27%  '
Set-TextValue $ws.Cells.Item(42, 4) '2.25'
$ws.Cells.Item(42, 5).Value = '  -1.77%  '
$ws.Cells.Item(43, 5).Value = '  +0.17%  '
Set-TextValue $ws.Cells.Item(44, 4) '1.789.09'
$ws.Cells.Item(44, 5).Value = '  +1.00%  '
$ws.Cells.Item(45, 5).Value = '  +0.83%  '
Set-TextValue $ws.Cells.Item(46, 4) '60.42'
$ws.Cells.Item(46, 5).Value = '  +9.87%  '
$ws.Cells.Item(47, 5).Value = '  +3.76%  '
$ws.Cells.Item(48, 5).Value = '  +0.62%  '
Set-TextValue $ws.Cells.Item(49, 4) '7.84'
Set-TextValue $ws.Cells.Item(50, 4) '0.0980'
$ws.Cells.Item(50, 5).Value = '  +3.22%  '
$ws.Cells.Item(51, 5).Value = '  -0.53%  '
